$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all cells we touch so Excel does not
# reinterpret numeric-looking or URL-looking strings as numbers/links.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.943.88'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.895.26'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7741'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.00'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3130'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07362'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08072'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7722'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.505'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.925.67'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.32'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.229'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.994.76'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.00'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.50'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007846'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.199.88'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.27%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.171'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1580'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.465'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.20'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.76'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.427'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.543'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.475'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.72%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.069'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.242'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7541'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.007'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.683'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01933'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.793'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.58'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4474'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.108.06'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +6.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.020'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8512'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.896'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.48'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.541'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.791'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.013'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.61%  '
